# ---------------------------------------------------------------------------
# platform schema, route, seed
#
# 1. Add a new "GetPlatform" worksheet (after "Change password") documenting
#    GET("/platforms") with its sample JSON response.
# 2. Move the active-tab selection from "Change password" (index 5) to
#    "register" (index 1).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsChangePassword = $wb.Worksheets.Item("Change password")
$wsRegister        = $wb.Worksheets.Item("register")

# --- 1. Update the selection left behind on "Change password" -------------
# (its own tab will no longer be the active one, and the previous
# topLeftCell scroll position goes away together with the selection change)
$wsChangePassword.Range("B4:E4").Select() | Out-Null

# --- 2. Create the new "GetPlatform" sheet, right after "Change password" -
$newSheet = $wb.Worksheets.Add($null, $wsChangePassword)
$newSheet.Name = "GetPlatform"

# Route header
$newSheet.Range("A2").Value = 'GET("/platforms")'

# Parameter-table header row, copied (with formatting) from an existing sheet
$wsChangePassword.Range("B4:E4").Copy($newSheet.Range("B4"))

# Response column header cells (style only, via copy), then set the real text
$wsChangePassword.Range("F10:G10").Copy($newSheet.Range("F7"))
$newSheet.Range("F7").Value = "Success"

$platformsJson = @'
{
    "status": 200,
    "data": [
        {
            "id": "6143fd832857c0ea5ecfa0e2",
            "name": "PS5"
        },
        {
            "id": "6143fd832857c0ea5ecfa0e3",
            "name": "PS4"
        },
        {
            "id": "6143fd832857c0ea5ecfa0e4",
            "name": "Xbox"
        },
        {
            "id": "6143fd832857c0ea5ecfa0e5",
            "name": "PC"
        },
        {
            "id": "6143fd832857c0ea5ecfa0e6",
            "name": "Nintendo Switch"
        }
    ],
    "error": null,
    "message": "Get platforms successfully"
}
'@
$newSheet.Range("G7").Value = $platformsJson

# Column widths / row height to match the other API sheets
$newSheet.Columns.Item(6).ColumnWidth = 8.14
$newSheet.Columns.Item(7).ColumnWidth = 40.43
$newSheet.Rows.Item(7).RowHeight = 372.6

$newSheet.Range("G4").Select() | Out-Null

# --- 3. Make "register" the active tab again -------------------------------
$wsRegister.Activate() | Out-Null
